$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.186.24'
$ws.Range("E2").Value = '  +0.87%  '
$ws.Range("D3").Value = '''1.851.56'
$ws.Range("E3").Value = '  +1.58%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = '''313.36'
$ws.Range("E5").Value = '  +0.69%  '
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("D7").Value = '''0.4643'
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("D8").Value = '''0.3712'
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("D9").Value = '''0.07279'
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("D10").Value = '''0.8862'
$ws.Range("E10").Value = '  +1.43%  '
$ws.Range("E11").Value = '  +1.62%  '
$ws.Range("D12").Value = '''0.07837'
$ws.Range("E12").Value = '  -0.68%  '
$ws.Range("D13").Value = '''1.797.42'
$ws.Range("E13").Value = '  -2.61%  '
$ws.Range("E14").Value = '  +0.91%  '
$ws.Range("D15").Value = '''6.519'
$ws.Range("E15").Value = '  -0.59%  '
$ws.Range("D16").Value = '''90.79'
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("E17").Value = '  -0.32%  '
$ws.Range("D18").Value = '''0.000008919'
$ws.Range("E18").Value = '  +0.98%  '
$ws.Range("E19").Value = '  -0.49%  '
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("D21").Value = '''27.221.38'
$ws.Range("E21").Value = '  +0.88%  '
$ws.Range("D22").Value = '''5.069'
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").Value = '''2.049.40'
$ws.Range("E24").Value = '  -2.89%  '
$ws.Range("D25").Value = '''1.953'
$ws.Range("E25").Value = '  +5.68%  '
$ws.Range("D26").Value = '''151.72'
$ws.Range("E26").Value = '  -1.15%  '
$ws.Range("D27").Value = '''18.35'
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("D28").Value = '''2.034'
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = '''115.71'
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").Value = '''5.049'
$ws.Range("E30").Value = '  -1.42%  '
$ws.Range("D31").Value = '''0.08801'
$ws.Range("E31").Value = '  -0.85%  '
$ws.Range("D32").Value = '''3.144'
$ws.Range("E32").Value = '  +6.20%  '
$ws.Range("D33").Value = '''0.7655'
$ws.Range("E33").Value = '  +5.08%  '
$ws.Range("E34").Value = '  +3.17%  '
$ws.Range("D35").Value = '''4.504'
$ws.Range("E35").Value = '  +1.55%  '
$ws.Range("D36").Value = '''2.711'
$ws.Range("E36").Value = '  +9.94%  '
$ws.Range("D37").Value = '''1.117'
$ws.Range("E37").Value = '  +4.53%  '
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("D39").Value = '''0.05202'
$ws.Range("E39").Value = '  -0.39%  '
$ws.Range("D40").Value = '''2.935'
$ws.Range("E40").Value = '  -0.42%  '
$ws.Range("D41").Value = '''7.013'
$ws.Range("E41").Value = '  -1.17%  '
$ws.Range("D42").Value = '''0.5099'
$ws.Range("E42").Value = '  -0.99%  '
$ws.Range("D43").Value = '''0.1627'
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").Value = '''8.410'
$ws.Range("E44").Value = '  +3.13%  '
$ws.Range("D45").Value = '''0.4790'
$ws.Range("E45").Value = '  -0.72%  '
$ws.Range("D46").Value = '''10.36'
$ws.Range("E46").Value = '  +1.94%  '
$ws.Range("D48").Value = '''102.79'
$ws.Range("E48").Value = '  +0.13%  '
$ws.Range("D49").Value = '''1.636'
$ws.Range("E49").Value = '  +0.39%  '
$ws.Range("D50").Value = '''0.06207'
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("D51").Value = '''65.54'
$ws.Range("E51").Value = '  +1.16%  '

# Reset style on the Price column so the quote-prefix trick used above
# does not leave a stray style index on any cell (keeps cells at the
# default/general style, matching the original workbook).
$ws.Range("D2:D51").Style = "Normal"
